# Add a new parameter row (motor_centerline_height) to the parameters table,
# matching the "Update STLs for 7 mm board thickness and 80 mm wheel" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "motor_centerline_height"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = "mm"
$ws.Range("D7").Value = "Height of motor shaft center from top of hardboard"

# Excel moves the active selection to the next empty row after data entry.
[void]$ws.Range("A8").Select()
